$wb = $excel.ActiveWorkbook

# Rename the second worksheet (sheetId=2, r:id=rId2) from "LoginTest" to "LoginFunc"
$wsLoginTest = $wb.Worksheets.Item("LoginTest")
$wsLoginTest.Name = "LoginFunc"

# Make the "LoginFunc" sheet the active sheet. This sets tabSelected on it,
# clears tabSelected on the previously active "devTestLogin" sheet, and
# updates the workbook's activeTab to this sheet's (0-based) index.
$wsLoginTest.Activate()

# Update the selection on the now-active sheet to D10.
$wsLoginTest.Range("D10").Select()
